# Atualização de bases das ligas, do dia: 10-06-2024 às 07:08
#
# For a handful of match rows, the row that should have held the "away" fixture
# data had been swapped with its neighbouring row (id, teams, odds, results, ...).
# This script restores the correct pairing by swapping back the contents of
# columns B (id) and E:AD (HomeTeam .. PL_AhUnder) between each affected pair of
# rows. Columns A (row index) and C/D (Div/Date) are identical between the two
# rows in every pair, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column helper: numeric column index -> used with Cells.Item(row, col)
$colB    = 2   # id
$colFrom = 5   # E : HomeTeam
$colTo   = 30  # AD: PL_AhUnder

function Swap-Cell($Row1, $Row2, $Col) {
    $cell1 = $ws.Cells.Item($Row1, $Col)
    $cell2 = $ws.Cells.Item($Row2, $Col)
    $val1 = $cell1.Value()
    $val2 = $cell2.Value()
    $cell1.Value = $val2
    $cell2.Value = $val1
}

function Swap-Rows($Row1, $Row2) {
    Swap-Cell $Row1 $Row2 $colB
    for ($col = $colFrom; $col -le $colTo; $col++) {
        Swap-Cell $Row1 $Row2 $col
    }
}

# Row pairs whose data needs to be swapped back
$pairs = @(
    @(31, 32),
    @(44, 45),
    @(59, 60),
    @(86, 87),
    @(164, 165),
    @(204, 205),
    @(222, 223),
    @(224, 225)
)

foreach ($pair in $pairs) {
    Swap-Rows $pair[0] $pair[1]
}
